$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: RRID becomes a numeric id, other cells keep the same text
$ws.Range("A2").Value = 22346

# --- Row 3: RRID becomes a numeric id, "Heat" parameter renamed to "Temperature"
$ws.Range("A3").Value = 12354
$ws.Range("F3").Value = "Temperature"

# --- Row 4 (new): soldering iron entry
$ws.Range("A4").Value = 65451
$ws.Range("B4").Value = "www.solderiron.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "http://www.solderiron.com/") | Out-Null
$ws.Range("B4").Style = $ws.Range("B2").Style
$ws.Range("C4").Value = "IronThrone"
$ws.Range("D4").Value = "2.4.6"
$ws.Range("E4").Value = "Wait until its hot"
$ws.Range("F4").Value = "Soldered"

# --- Row 5 (new): microphone entry, placed in front
$ws.Range("A5").Value = 32142
$ws.Range("B5").Value = "www.Microphone.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "http://www.Microphone.com/") | Out-Null
$ws.Range("B5").Style = $ws.Range("B2").Style
$ws.Range("C5").Value = "Hot mics"
$ws.Range("D5").Value = "2.2.2"
$ws.Range("E5").Value = "Placed in front of mouse"
$ws.Range("F5").Value = "Sound"
$ws.Range("F5").Font.Name = "Calibri"

# --- Row 6 (new): microphone entry, placed behind
$ws.Range("A6").Value = 35543
$ws.Range("B6").Value = "www.Microphone.com"
$ws.Hyperlinks.Add($ws.Range("B6"), "http://www.Microphone.com/") | Out-Null
$ws.Range("B6").Style = $ws.Range("B2").Style
$ws.Range("C6").Value = "Hot mics"
$ws.Range("D6").Value = "1.1.4"
$ws.Range("E6").Value = "Placed behind mouse"
$ws.Range("F6").Value = "Sound"
$ws.Range("F6").Font.Name = "Calibri"

# --- Row 7 (new): microphone entry, placed above
$ws.Range("A7").Value = 33242
$ws.Range("B7").Value = "www.Microphone.com"
$ws.Hyperlinks.Add($ws.Range("B7"), "http://www.Microphone.com/") | Out-Null
$ws.Range("B7").Style = $ws.Range("B2").Style
$ws.Range("C7").Value = "Hot mics"
$ws.Range("D7").Value = "4.4.9"
$ws.Range("E7").Value = "Placed above mouse"
$ws.Range("F7").Value = "Sound"
$ws.Range("F7").Font.Name = "Calibri"

# --- Selection moves to F8, mirroring the saved cursor position in the source file
$ws.Range("F8").Select() | Out-Null
